$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the daily conversion text in A1 ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 3.14 = 11772.05 pesos`n✅ 11772.05 pesos = 3.12 = 970.42 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- Sheet "tasas": update the rate cells ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 317.999
$ws2.Range("O10").Value = 3743.5
$ws2.Range("N12").Value = 3772.7
